$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-29: 45243 -> 45244
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
